$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two brand-new columns at B:C; this shifts the existing
# "Jun_13" column (old B) to D and the existing "Jun_10" column
# (old C) to E -- each cell's style moves along with its value.
$ws.Range("B1:C1").EntireColumn.Insert()

# New column headers: B gets the newest watch date, C the next newest.
$ws.Range("B1").Value = "Jun_17"
$ws.Range("C1").Value = "Jun_15"

# Fill the two freshly inserted data columns with the same placeholder
# rating ("UN") used throughout the rest of the sheet.
$ws.Range("B2:C27").Value = "UN"

# Match the column widths used elsewhere on the sheet (stored width 8.0).
$ws.Range("C1:E1").EntireColumn.ColumnWidth = 7.14
